$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.046.80"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.530.59"
$ws.Range("E3").Value = "  -1.12%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.02"
$ws.Range("E5").Value = "  -1.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.38"
$ws.Range("E6").Value = "  -2.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.526.19"
$ws.Range("E7").Value = "  -1.19%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.493"
$ws.Range("E9").Value = "  -1.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.125"
$ws.Range("E10").Value = "  +0.66%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.15"
$ws.Range("E11").Value = "  +2.81%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.387"
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.127.72"
$ws.Range("E13").Value = "  -1.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.72"
$ws.Range("E14").Value = "  +1.13%  "
$ws.Range("E15").Value = "  -1.07%  "
$ws.Range("E16").Value = "  +0.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.528.15"
$ws.Range("E17").Value = "  -1.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "65.029.67"
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.10"
$ws.Range("E19").Value = "  -0.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.42"
$ws.Range("E20").Value = "  -0.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.72"
$ws.Range("E21").Value = "  -2.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "392.36"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.669.20"
$ws.Range("E24").Value = "  -1.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.56"
$ws.Range("E25").Value = "  +0.49%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000113"
$ws.Range("E27").Value = "  -3.95%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.72"
$ws.Range("E28").Value = "  +0.63%  "
$ws.Range("B29").Value = "Fetch.AI"
$ws.Range("C29").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.57"
$ws.Range("E29").Value = "  +8.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -0.31%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.39"
$ws.Range("E31").Value = "  +0.38%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.27"
$ws.Range("E32").Value = "  -0.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.530.31"
$ws.Range("E33").Value = "  -1.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "24.24"
$ws.Range("E34").Value = "  +0.63%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  +0.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.28"
$ws.Range("E37").Value = "  +5.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.02"
$ws.Range("E38").Value = "  +0.91%  "
$ws.Range("E39").Value = "  +0.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "168.77"
$ws.Range("E40").Value = "  -1.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0818"
$ws.Range("E41").Value = "  +0.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.825"
$ws.Range("E42").Value = "  -0.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.27"
$ws.Range("E43").Value = "  +3.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "25.76"
$ws.Range("E44").Value = "  -5.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.89"
$ws.Range("E45").Value = "  +0.44%  "
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("E47").Value = "  -1.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.66"
$ws.Range("E48").Value = "  -0.78%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.94"
$ws.Range("E49").Value = "  +0.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.416.10"
$ws.Range("E50").Value = "  -2.55%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.910"
$ws.Range("E51").Value = "  +5.44%  "
